# Sync the GitHub copy of this workbook with the current state of the
# author's main working copy: drop the two now-unused scratch sheets and
# refresh Sheet1's computed values; also leave A1:D4 selected, matching
# the last on-screen state in the source workbook.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet2 and Sheet3 were scratch/unused sheets - remove them.
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# Refresh the computed sensorCG X/Y/Z values on Sheet1 (row 1 headers,
# column A unchanged).
$ws.Range("B2").Value = 0.97006622951441601
$ws.Range("C2").Value = 0.17280460603734368
$ws.Range("D2").Value = -0.17061675910637428

$ws.Range("B3").Value = 0.23992437651154333
$ws.Range("C3").Value = -0.57345181507739651
$ws.Range("D3").Value = 0.78331941718559916

$ws.Range("B4").Value = 0.037520713095985883
$ws.Range("C4").Value = -0.80080683308568101
$ws.Range("D4").Value = -0.5977462774221608

$ws.Range("A1:D4").Select()
